$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(2, 6).Value = 2538
$ws.Cells.Item(3, 6).Value = 548
$ws.Cells.Item(4, 6).Value = 450
$ws.Cells.Item(5, 6).Value = 283
$ws.Cells.Item(8, 6).Value = 1182
$ws.Cells.Item(9, 6).Value = 530
$ws.Cells.Item(10, 6).Value = 290
$ws.Cells.Item(11, 6).Value = 111
$ws.Cells.Item(13, 6).Value = 5472
$ws.Cells.Item(15, 6).Value = 1653
$ws.Cells.Item(16, 6).Value = 3973
$ws.Cells.Item(17, 6).Value = 399
$ws.Cells.Item(20, 6).Value = 4520
$ws.Cells.Item(21, 6).Value = 5914
$ws.Cells.Item(24, 6).Value = 649
$ws.Cells.Item(25, 6).Value = 3657
$ws.Cells.Item(26, 6).Value = 462
$ws.Cells.Item(30, 6).Value = 958
$ws.Cells.Item(31, 6).Value = 1353
$ws.Cells.Item(32, 6).Value = 443
$ws.Cells.Item(33, 6).Value = 504
$ws.Cells.Item(34, 6).Value = 1547
$ws.Cells.Item(35, 6).Value = 186
$ws.Cells.Item(36, 6).Value = 1625
$ws.Cells.Item(37, 6).Value = 152
$ws.Cells.Item(38, 6).Value = 2
$ws.Cells.Item(39, 6).Value = 1067
$ws.Cells.Item(40, 6).Value = 29
$ws.Cells.Item(41, 6).Value = 1341
$ws.Cells.Item(42, 6).Value = 597
$ws.Cells.Item(44, 6).Value = 204
$ws.Cells.Item(44, 7).Value = "不可售"
$ws.Cells.Item(45, 6).Value = 2748
$ws.Cells.Item(46, 6).Value = 118
$ws.Cells.Item(47, 6).Value = 256
$ws.Cells.Item(49, 6).Value = 3852
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(5, 6).Value = 1161
$ws.Cells.Item(20, 6).Value = 52
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(2, 6).Value = 3647
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(3, 6).Value = 2538
$ws.Cells.Item(4, 6).Value = 548
$ws.Cells.Item(5, 6).Value = 450
$ws.Cells.Item(6, 6).Value = 283
$ws.Cells.Item(7, 6).Value = 1161
$ws.Cells.Item(10, 6).Value = 1183
$ws.Cells.Item(11, 6).Value = 530
$ws.Cells.Item(12, 6).Value = 290
$ws.Cells.Item(13, 6).Value = 111
$ws.Cells.Item(15, 6).Value = 5472
$ws.Cells.Item(17, 6).Value = 1653
$ws.Cells.Item(18, 6).Value = 4521
$ws.Cells.Item(19, 6).Value = 5914
$ws.Cells.Item(22, 6).Value = 649
$ws.Cells.Item(23, 6).Value = 3657
$ws.Cells.Item(24, 6).Value = 462
$ws.Cells.Item(28, 6).Value = 1353
$ws.Cells.Item(29, 6).Value = 443
$ws.Cells.Item(30, 6).Value = 504
$ws.Cells.Item(32, 6).Value = 1547
$ws.Cells.Item(33, 6).Value = 186
$ws.Cells.Item(34, 6).Value = 1625
$ws.Cells.Item(36, 6).Value = 1067
$ws.Cells.Item(38, 6).Value = 597
$ws.Cells.Item(39, 6).Value = 52
$ws.Cells.Item(43, 6).Value = 2748
$ws.Cells.Item(45, 6).Value = 118
$ws.Cells.Item(46, 6).Value = 256
$ws.Cells.Item(49, 6).Value = 3852
